$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.692.81"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "3.500.79"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.82%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.498.70"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").Value = "4.105.88"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "66.710.11"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "3.496.60"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.533"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("E31").Value = "  +5.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0747"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.821.69"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "356.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.34%  "
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.89%  "
